$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coefficient values in column B
$ws.Range("B2").Value = 0.5410592318281237
$ws.Range("B3").Value = -0.34902016125534
$ws.Range("B4").Value = 101.36095924155676

# Remove row 5 (label "4" and its coefficient) entirely
$ws.Range("A5:B5").Delete()
